$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 51.310475741437763
$ws.Range("C3").Value = 50.170464970209075

$ws.Range("B4").Value = 51.752756646435522
$ws.Range("C4").Value = 50.187732508154959

$ws.Range("B5").Value = 52.488874041344161
$ws.Range("C5").Value = 50.659118779642291

$ws.Range("B6").Value = 52.684871806895259
$ws.Range("C6").Value = 50.695543352224973

$ws.Range("B7").Value = 53.152730033410855
$ws.Range("C7").Value = 51.062618369904165

$ws.Range("B9").Value = -0.34470847390136639
$ws.Range("C9").Value = -0.36121802038602074

$ws.Range("B13").Value = 0.81769406092511432
$ws.Range("C13").Value = 0.96726330638930747

$ws.Range("B14").Value = 41.941244041736169
$ws.Range("C14").Value = 49.059363041775271

$ws.Range("B15").Value = 51.516132699057621
$ws.Range("C15").Value = 51.110757709397525

$ws.Range("C17").Value = -0.2291458696175728

$ws.Range("C18").Value = -0.10413878312056032

$ws.Range("C19").Value = -0.087849378515865245

$ws.Range("C20").Value = -1.4345073382850331

$ws.Range("C21").Value = -0.043867024355286756

$ws.Range("B22").Value = 8050
$ws.Range("C22").Value = 9932
